# Applies the cryptos-list refresh described in the commit message.
# For numeric-looking Price (column D) values we pin the cell to Text
# first (otherwise assigning e.g. 9.00 would silently become the number
# 9) and restore the default "Normal" style afterwards so the only
# observable change is the cell's value, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.389.42"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "3.673.13"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "639.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +0.50%  "

$ws.Range("E9").Value = "  -1.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.445"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000231"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").Value = "4.291.60"
$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.04%  "

$ws.Range("D15").Value = "3.684.35"
$ws.Range("E15").Value = "  -0.31%  "

$ws.Range("D16").Value = "69.364.34"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("E22").Value = "  -0.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.16%  "

$ws.Range("D24").Value = "3.819.74"
$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000125"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.29%  "

$ws.Range("E30").Value = "  -1.92%  "

$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("E32").Value = "  -0.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.76%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.163"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.36%  "

$ws.Range("D36").Value = "3.666.09"
$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.61%  "

$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.38%  "

$ws.Range("E42").Value = "  -0.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0897"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.924"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.06%  "

$ws.Range("B48").Value = "SuiNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.28%  "

$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000268"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.33%  "

